$wb = $excel.ActiveWorkbook

# --- Sheet1: Recommandations ---
$ws1 = $wb.Worksheets.Item("Recommandations")

# Row 2: NEI-CEDA CI
$ws1.Cells.Item(2, 1).Value = 'NEI-CEDA CI'
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 4
$ws1.Cells.Item(2, 4).Value = 3700
$ws1.Cells.Item(2, 5).Value = 930
$ws1.Cells.Item(2, 6).Value = '🟡 Observer'
$ws1.Cells.Item(2, 7).Value = '➖ Neutre'

# Row 3: BRVM - SERVICES PUBLICS
$ws1.Cells.Item(3, 1).Value = 'BRVM - SERVICES PUBLICS'
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 8
$ws1.Cells.Item(3, 4).Value = 3350.44
$ws1.Cells.Item(3, 5).Value = 111.67
$ws1.Cells.Item(3, 6).Value = '🟡 Observer'
$ws1.Cells.Item(3, 7).Value = '➖ Neutre'

# Row 4: BRVM - AUTRES SECTEURS
$ws1.Cells.Item(4, 1).Value = 'BRVM - AUTRES SECTEURS'
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 4
$ws1.Cells.Item(4, 4).Value = 2364.43
$ws1.Cells.Item(4, 5).Value = 596.54
$ws1.Cells.Item(4, 6).Value = '🟡 Observer'
$ws1.Cells.Item(4, 7).Value = '➖ Neutre'

# Row 5: BRVM - DISTRIBUTION
$ws1.Cells.Item(5, 1).Value = 'BRVM - DISTRIBUTION'
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 4
$ws1.Cells.Item(5, 4).Value = 2011.08
$ws1.Cells.Item(5, 5).Value = 502.54
$ws1.Cells.Item(5, 6).Value = '🟡 Observer'
$ws1.Cells.Item(5, 7).Value = '➖ Neutre'

# Row 6: BRVM - TRANSPORT
$ws1.Cells.Item(6, 1).Value = 'BRVM - TRANSPORT'
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 4
$ws1.Cells.Item(6, 4).Value = 1434.37
$ws1.Cells.Item(6, 5).Value = 364.71
$ws1.Cells.Item(6, 6).Value = '🟡 Observer'
$ws1.Cells.Item(6, 7).Value = '➖ Neutre'

# Row 7: BRVM - AGRICULTURE
$ws1.Cells.Item(7, 1).Value = 'BRVM - AGRICULTURE'
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 4
$ws1.Cells.Item(7, 4).Value = 1351.74
$ws1.Cells.Item(7, 5).Value = 336.59
$ws1.Cells.Item(7, 6).Value = '🟡 Observer'
$ws1.Cells.Item(7, 7).Value = '➖ Neutre'

# Row 8: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Cells.Item(8, 1).Value = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 4
$ws1.Cells.Item(8, 4).Value = 680.52
$ws1.Cells.Item(8, 5).Value = 171.1
$ws1.Cells.Item(8, 6).Value = '🟡 Observer'
$ws1.Cells.Item(8, 7).Value = '➖ Neutre'

# Row 9: BRVM - FINANCES
$ws1.Cells.Item(9, 1).Value = 'BRVM - FINANCES'
$ws1.Cells.Item(9, 2).Value = 0
$ws1.Cells.Item(9, 3).Value = 4
$ws1.Cells.Item(9, 4).Value = 588.48
$ws1.Cells.Item(9, 5).Value = 146.37
$ws1.Cells.Item(9, 6).Value = '🟡 Observer'
$ws1.Cells.Item(9, 7).Value = '➖ Neutre'

# Row 10: BRVM - SERVICES FINANCIERS
$ws1.Cells.Item(10, 1).Value = 'BRVM - SERVICES FINANCIERS'
$ws1.Cells.Item(10, 2).Value = 0
$ws1.Cells.Item(10, 3).Value = 4
$ws1.Cells.Item(10, 4).Value = 578.36
$ws1.Cells.Item(10, 5).Value = 143.85
$ws1.Cells.Item(10, 6).Value = '🟡 Observer'
$ws1.Cells.Item(10, 7).Value = '➖ Neutre'

# Row 11: BRVM-PRESTIGE
$ws1.Cells.Item(11, 1).Value = 'BRVM-PRESTIGE'
$ws1.Cells.Item(11, 2).Value = 0
$ws1.Cells.Item(11, 3).Value = 4
$ws1.Cells.Item(11, 4).Value = 567.58
$ws1.Cells.Item(11, 5).Value = 140.51
$ws1.Cells.Item(11, 6).Value = '🟡 Observer'
$ws1.Cells.Item(11, 7).Value = '➖ Neutre'

# Row 12: BRVM - INDUSTRIE                 (**)
$ws1.Cells.Item(12, 1).Value = 'BRVM - INDUSTRIE                 (**)'
$ws1.Cells.Item(12, 2).Value = 0
$ws1.Cells.Item(12, 3).Value = 2
$ws1.Cells.Item(12, 4).Value = 533.8200000000001
$ws1.Cells.Item(12, 5).Value = 266.86
$ws1.Cells.Item(12, 6).Value = '🟡 Observer'
$ws1.Cells.Item(12, 7).Value = '➖ Neutre'

# Row 13: BRVM - INDUSTRIELS
$ws1.Cells.Item(13, 1).Value = 'BRVM - INDUSTRIELS'
$ws1.Cells.Item(13, 2).Value = 0
$ws1.Cells.Item(13, 3).Value = 4
$ws1.Cells.Item(13, 4).Value = 505.94
$ws1.Cells.Item(13, 5).Value = 128.63
$ws1.Cells.Item(13, 6).Value = '🟡 Observer'
$ws1.Cells.Item(13, 7).Value = '➖ Neutre'

# Row 14: BRVM - ENERGIE
$ws1.Cells.Item(14, 1).Value = 'BRVM - ENERGIE'
$ws1.Cells.Item(14, 2).Value = 0
$ws1.Cells.Item(14, 3).Value = 4
$ws1.Cells.Item(14, 4).Value = 451.76
$ws1.Cells.Item(14, 5).Value = 111.97
$ws1.Cells.Item(14, 6).Value = '🟡 Observer'
$ws1.Cells.Item(14, 7).Value = '➖ Neutre'

# Row 15: BRVM - TELECOMMUNICATIONS
$ws1.Cells.Item(15, 1).Value = 'BRVM - TELECOMMUNICATIONS'
$ws1.Cells.Item(15, 2).Value = 0
$ws1.Cells.Item(15, 3).Value = 4
$ws1.Cells.Item(15, 4).Value = 377.95
$ws1.Cells.Item(15, 5).Value = 92.92
$ws1.Cells.Item(15, 6).Value = '🟡 Observer'
$ws1.Cells.Item(15, 7).Value = '➖ Neutre'

# Row 16: BRVM - INDUSTRIE
$ws1.Cells.Item(16, 1).Value = 'BRVM - INDUSTRIE'
$ws1.Cells.Item(16, 2).Value = 0
$ws1.Cells.Item(16, 3).Value = 1
$ws1.Cells.Item(16, 4).Value = 266.4
$ws1.Cells.Item(16, 5).Value = 266.4
$ws1.Cells.Item(16, 6).Value = '🟡 Observer'
$ws1.Cells.Item(16, 7).Value = '➖ Neutre'

# Row 17: BRVM - CONSOMMATION DE BASE          (**)
$ws1.Cells.Item(17, 1).Value = 'BRVM - CONSOMMATION DE BASE          (**)'
$ws1.Cells.Item(17, 2).Value = 0
$ws1.Cells.Item(17, 3).Value = 1
$ws1.Cells.Item(17, 4).Value = 222.9
$ws1.Cells.Item(17, 5).Value = 222.9
$ws1.Cells.Item(17, 6).Value = '🟡 Observer'
$ws1.Cells.Item(17, 7).Value = '➖ Neutre'

# Row 18: BRVM - CONSOMMATION DE BASE         (**)
$ws1.Cells.Item(18, 1).Value = 'BRVM - CONSOMMATION DE BASE         (**)'
$ws1.Cells.Item(18, 2).Value = 0
$ws1.Cells.Item(18, 3).Value = 1
$ws1.Cells.Item(18, 4).Value = 222.7
$ws1.Cells.Item(18, 5).Value = 222.7
$ws1.Cells.Item(18, 6).Value = '🟡 Observer'
$ws1.Cells.Item(18, 7).Value = '➖ Neutre'

# Row 19: BRVM - CONSOMMATION DE BASE
$ws1.Cells.Item(19, 1).Value = 'BRVM - CONSOMMATION DE BASE'
$ws1.Cells.Item(19, 2).Value = 0
$ws1.Cells.Item(19, 3).Value = 1
$ws1.Cells.Item(19, 4).Value = 222.06
$ws1.Cells.Item(19, 5).Value = 222.06
$ws1.Cells.Item(19, 6).Value = '🟡 Observer'
$ws1.Cells.Item(19, 7).Value = '➖ Neutre'

# Row 20: BRVM-PRINCIPAL                    (**)
$ws1.Cells.Item(20, 1).Value = 'BRVM-PRINCIPAL                    (**)'
$ws1.Cells.Item(20, 2).Value = 0
$ws1.Cells.Item(20, 3).Value = 1
$ws1.Cells.Item(20, 4).Value = 220.55
$ws1.Cells.Item(20, 5).Value = 220.55
$ws1.Cells.Item(20, 6).Value = '🟡 Observer'
$ws1.Cells.Item(20, 7).Value = '➖ Neutre'

# Row 21: BRVM-PRINCIPAL                   (**)
$ws1.Cells.Item(21, 1).Value = 'BRVM-PRINCIPAL                   (**)'
$ws1.Cells.Item(21, 2).Value = 0
$ws1.Cells.Item(21, 3).Value = 1
$ws1.Cells.Item(21, 4).Value = 220.21
$ws1.Cells.Item(21, 5).Value = 220.21
$ws1.Cells.Item(21, 6).Value = '🟡 Observer'
$ws1.Cells.Item(21, 7).Value = '➖ Neutre'

# Row 22: BRVM-PRINCIPAL
$ws1.Cells.Item(22, 1).Value = 'BRVM-PRINCIPAL'
$ws1.Cells.Item(22, 2).Value = 0
$ws1.Cells.Item(22, 3).Value = 1
$ws1.Cells.Item(22, 4).Value = 220.02
$ws1.Cells.Item(22, 5).Value = 220.02
$ws1.Cells.Item(22, 6).Value = '🟡 Observer'
$ws1.Cells.Item(22, 7).Value = '➖ Neutre'

# Row 23: ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)
$ws1.Cells.Item(23, 1).Value = 'ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)'
$ws1.Cells.Item(23, 2).Value = 4
$ws1.Cells.Item(23, 3).Value = 0
$ws1.Cells.Item(23, 4).Value = 29.46
$ws1.Cells.Item(23, 5).Value = 7.49
$ws1.Cells.Item(23, 6).Value = '🟢 Achat'
$ws1.Cells.Item(23, 7).Value = '✅ Renforcer'

# Row 24: SERVAIR ABIDJAN CI (ABJC)
$ws1.Cells.Item(24, 1).Value = 'SERVAIR ABIDJAN CI (ABJC)'
$ws1.Cells.Item(24, 2).Value = 2
$ws1.Cells.Item(24, 3).Value = 0
$ws1.Cells.Item(24, 4).Value = 11.29
$ws1.Cells.Item(24, 5).Value = 3.92
$ws1.Cells.Item(24, 6).Value = '🟡 Observer'
$ws1.Cells.Item(24, 7).Value = '➖ Neutre'

# Row 25: SUCRIVOIRE (SCRC)
$ws1.Cells.Item(25, 1).Value = 'SUCRIVOIRE (SCRC)'
$ws1.Cells.Item(25, 2).Value = 2
$ws1.Cells.Item(25, 3).Value = 0
$ws1.Cells.Item(25, 4).Value = 9.359999999999999
$ws1.Cells.Item(25, 5).Value = 1.86
$ws1.Cells.Item(25, 6).Value = '🟡 Observer'
$ws1.Cells.Item(25, 7).Value = '➖ Neutre'

# Row 26: SOLIBRA CI (SLBC)
$ws1.Cells.Item(26, 1).Value = 'SOLIBRA CI (SLBC)'
$ws1.Cells.Item(26, 2).Value = 1
$ws1.Cells.Item(26, 3).Value = 0
$ws1.Cells.Item(26, 4).Value = 7.49
$ws1.Cells.Item(26, 5).Value = 7.49
$ws1.Cells.Item(26, 6).Value = '🟡 Observer'
$ws1.Cells.Item(26, 7).Value = '➖ Neutre'

# Row 27: BERNABE CI (BNBC)
$ws1.Cells.Item(27, 1).Value = 'BERNABE CI (BNBC)'
$ws1.Cells.Item(27, 2).Value = 2
$ws1.Cells.Item(27, 3).Value = 1
$ws1.Cells.Item(27, 4).Value = 6.45
$ws1.Cells.Item(27, 5).Value = 3.21
$ws1.Cells.Item(27, 6).Value = '🟡 Observer'
$ws1.Cells.Item(27, 7).Value = '👀 À surveiller'

# Row 28: ECOBANK TRANS. INCORP. TG (ETIT)
$ws1.Cells.Item(28, 1).Value = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$ws1.Cells.Item(28, 2).Value = 1
$ws1.Cells.Item(28, 3).Value = 0
$ws1.Cells.Item(28, 4).Value = 4.55
$ws1.Cells.Item(28, 5).Value = 4.55
$ws1.Cells.Item(28, 6).Value = '🟡 Observer'
$ws1.Cells.Item(28, 7).Value = '➖ Neutre'

# Row 29: BANK OF AFRICA NG (BOAN)
$ws1.Cells.Item(29, 1).Value = 'BANK OF AFRICA NG (BOAN)'
$ws1.Cells.Item(29, 2).Value = 1
$ws1.Cells.Item(29, 3).Value = 0
$ws1.Cells.Item(29, 4).Value = 3.56
$ws1.Cells.Item(29, 5).Value = 3.56
$ws1.Cells.Item(29, 6).Value = '🟡 Observer'
$ws1.Cells.Item(29, 7).Value = '➖ Neutre'

# Row 30: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Cells.Item(30, 1).Value = 'AFRICA GLOBAL LOGISTICS CI (SDSC)'
$ws1.Cells.Item(30, 2).Value = 1
$ws1.Cells.Item(30, 3).Value = 0
$ws1.Cells.Item(30, 4).Value = 2.76
$ws1.Cells.Item(30, 5).Value = 2.76
$ws1.Cells.Item(30, 6).Value = '🟡 Observer'
$ws1.Cells.Item(30, 7).Value = '➖ Neutre'

# Row 31: CORIS BANK INTERNATIONAL (CBIBF)
$ws1.Cells.Item(31, 1).Value = 'CORIS BANK INTERNATIONAL (CBIBF)'
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 3).Value = 0
$ws1.Cells.Item(31, 4).Value = 1.6
$ws1.Cells.Item(31, 5).Value = 1.6
$ws1.Cells.Item(31, 6).Value = '🟡 Observer'
$ws1.Cells.Item(31, 7).Value = '➖ Neutre'

# Row 32: FILTISAC CI (FTSC)
$ws1.Cells.Item(32, 1).Value = 'FILTISAC CI (FTSC)'
$ws1.Cells.Item(32, 2).Value = 1
$ws1.Cells.Item(32, 3).Value = 1
$ws1.Cells.Item(32, 4).Value = 0.08
$ws1.Cells.Item(32, 5).Value = 2.91
$ws1.Cells.Item(32, 6).Value = '🟡 Observer'
$ws1.Cells.Item(32, 7).Value = '👀 À surveiller'

# Row 33: TOTAL
$ws1.Cells.Item(33, 1).Value = 'TOTAL'
$ws1.Cells.Item(33, 2).Value = 0
$ws1.Cells.Item(33, 3).Value = 4
$ws1.Cells.Item(33, 4).Value = 0
$ws1.Cells.Item(33, 5).Value = 0
$ws1.Cells.Item(33, 6).Value = '🟡 Observer'
$ws1.Cells.Item(33, 7).Value = '➖ Neutre'

# Row 34: ONATEL BF (ONTBF)
$ws1.Cells.Item(34, 1).Value = 'ONATEL BF (ONTBF)'
$ws1.Cells.Item(34, 2).Value = 2
$ws1.Cells.Item(34, 3).Value = 1
$ws1.Cells.Item(34, 4).Value = -0.1
$ws1.Cells.Item(34, 5).Value = 2.04
$ws1.Cells.Item(34, 6).Value = '🟡 Observer'
$ws1.Cells.Item(34, 7).Value = '👀 À surveiller'

# Row 35: SICOR CI (SICC)
$ws1.Cells.Item(35, 1).Value = 'SICOR CI (SICC)'
$ws1.Cells.Item(35, 2).Value = 1
$ws1.Cells.Item(35, 3).Value = 2
$ws1.Cells.Item(35, 4).Value = -1.02
$ws1.Cells.Item(35, 5).Value = 3.85
$ws1.Cells.Item(35, 6).Value = '🟡 Observer'
$ws1.Cells.Item(35, 7).Value = '👀 À surveiller'

# Row 36: SMB CI (SMBC)
$ws1.Cells.Item(36, 1).Value = 'SMB CI (SMBC)'
$ws1.Cells.Item(36, 2).Value = 0
$ws1.Cells.Item(36, 3).Value = 1
$ws1.Cells.Item(36, 4).Value = -1.05
$ws1.Cells.Item(36, 5).Value = -1.05
$ws1.Cells.Item(36, 6).Value = '🟡 Observer'
$ws1.Cells.Item(36, 7).Value = '➖ Neutre'

# Row 37: BICI CI (BICC)
$ws1.Cells.Item(37, 1).Value = 'BICI CI (BICC)'
$ws1.Cells.Item(37, 2).Value = 0
$ws1.Cells.Item(37, 3).Value = 1
$ws1.Cells.Item(37, 4).Value = -1.64
$ws1.Cells.Item(37, 5).Value = -1.64
$ws1.Cells.Item(37, 6).Value = '🟡 Observer'
$ws1.Cells.Item(37, 7).Value = '➖ Neutre'

# Row 38: BANK OF AFRICA SENEGAL (BOAS)
$ws1.Cells.Item(38, 1).Value = 'BANK OF AFRICA SENEGAL (BOAS)'
$ws1.Cells.Item(38, 2).Value = 0
$ws1.Cells.Item(38, 3).Value = 1
$ws1.Cells.Item(38, 4).Value = -1.76
$ws1.Cells.Item(38, 5).Value = -1.76
$ws1.Cells.Item(38, 6).Value = '🟡 Observer'
$ws1.Cells.Item(38, 7).Value = '➖ Neutre'

# Row 39: SETAO CI (STAC)
$ws1.Cells.Item(39, 1).Value = 'SETAO CI (STAC)'
$ws1.Cells.Item(39, 2).Value = 0
$ws1.Cells.Item(39, 3).Value = 1
$ws1.Cells.Item(39, 4).Value = -1.87
$ws1.Cells.Item(39, 5).Value = -1.87
$ws1.Cells.Item(39, 6).Value = '🟡 Observer'
$ws1.Cells.Item(39, 7).Value = '➖ Neutre'

# Row 40: BANK OF AFRICA CI (BOAC)
$ws1.Cells.Item(40, 1).Value = 'BANK OF AFRICA CI (BOAC)'
$ws1.Cells.Item(40, 2).Value = 0
$ws1.Cells.Item(40, 3).Value = 1
$ws1.Cells.Item(40, 4).Value = -2.3
$ws1.Cells.Item(40, 5).Value = -2.3
$ws1.Cells.Item(40, 6).Value = '🟡 Observer'
$ws1.Cells.Item(40, 7).Value = '➖ Neutre'

# Row 41: SOCIETE IVOIRIENNE DE BANQUE  (SIBC)
$ws1.Cells.Item(41, 1).Value = 'SOCIETE IVOIRIENNE DE BANQUE  (SIBC)'
$ws1.Cells.Item(41, 2).Value = 0
$ws1.Cells.Item(41, 3).Value = 1
$ws1.Cells.Item(41, 4).Value = -2.42
$ws1.Cells.Item(41, 5).Value = -2.42
$ws1.Cells.Item(41, 6).Value = '🟡 Observer'
$ws1.Cells.Item(41, 7).Value = '➖ Neutre'

# Row 42: ORAGROUP TOGO (ORGT)
$ws1.Cells.Item(42, 1).Value = 'ORAGROUP TOGO (ORGT)'
$ws1.Cells.Item(42, 2).Value = 0
$ws1.Cells.Item(42, 3).Value = 1
$ws1.Cells.Item(42, 4).Value = -2.72
$ws1.Cells.Item(42, 5).Value = -2.72
$ws1.Cells.Item(42, 6).Value = '🟡 Observer'
$ws1.Cells.Item(42, 7).Value = '➖ Neutre'

# Row 43: NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws1.Cells.Item(43, 1).Value = 'NSIA BANQUE COTE D''IVOIRE (NSBC)'
$ws1.Cells.Item(43, 2).Value = 0
$ws1.Cells.Item(43, 3).Value = 1
$ws1.Cells.Item(43, 4).Value = -3.48
$ws1.Cells.Item(43, 5).Value = -3.48
$ws1.Cells.Item(43, 6).Value = '🟡 Observer'
$ws1.Cells.Item(43, 7).Value = '➖ Neutre'

# Row 44: VIVO ENERGY CI (SHEC)
$ws1.Cells.Item(44, 1).Value = 'VIVO ENERGY CI (SHEC)'
$ws1.Cells.Item(44, 2).Value = 0
$ws1.Cells.Item(44, 3).Value = 1
$ws1.Cells.Item(44, 4).Value = -3.85
$ws1.Cells.Item(44, 5).Value = -3.85
$ws1.Cells.Item(44, 6).Value = '🟡 Observer'
$ws1.Cells.Item(44, 7).Value = '➖ Neutre'

# Row 45: SICABLE CI (CABC)
$ws1.Cells.Item(45, 1).Value = 'SICABLE CI (CABC)'
$ws1.Cells.Item(45, 2).Value = 0
$ws1.Cells.Item(45, 3).Value = 1
$ws1.Cells.Item(45, 4).Value = -3.89
$ws1.Cells.Item(45, 5).Value = -3.89
$ws1.Cells.Item(45, 6).Value = '🟡 Observer'
$ws1.Cells.Item(45, 7).Value = '➖ Neutre'

# Row 46: NEI-CEDA CI (NEIC)
$ws1.Cells.Item(46, 1).Value = 'NEI-CEDA CI (NEIC)'
$ws1.Cells.Item(46, 2).Value = 0
$ws1.Cells.Item(46, 3).Value = 1
$ws1.Cells.Item(46, 4).Value = -4.26
$ws1.Cells.Item(46, 5).Value = -4.26
$ws1.Cells.Item(46, 6).Value = '🟡 Observer'
$ws1.Cells.Item(46, 7).Value = '➖ Neutre'

# Row 47: ORANGE COTE D'IVOIRE (ORAC)
$ws1.Cells.Item(47, 1).Value = 'ORANGE COTE D''IVOIRE (ORAC)'
$ws1.Cells.Item(47, 2).Value = 0
$ws1.Cells.Item(47, 3).Value = 2
$ws1.Cells.Item(47, 4).Value = -5.98
$ws1.Cells.Item(47, 5).Value = -4.03
$ws1.Cells.Item(47, 6).Value = '🟡 Observer'
$ws1.Cells.Item(47, 7).Value = '➖ Neutre'

# Row 48: UNILEVER CI (UNLC)
$ws1.Cells.Item(48, 1).Value = 'UNILEVER CI (UNLC)'
$ws1.Cells.Item(48, 2).Value = 0
$ws1.Cells.Item(48, 3).Value = 1
$ws1.Cells.Item(48, 4).Value = -6.97
$ws1.Cells.Item(48, 5).Value = -6.97
$ws1.Cells.Item(48, 6).Value = '🟡 Observer'
$ws1.Cells.Item(48, 7).Value = '➖ Neutre'

# --- Sheet2: Top_YTD ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

# Row 2: BRVM - SERVICES PUBLICS
$ws2.Cells.Item(2, 2).Value = 9321255.16

# Row 3: NEI-CEDA CI
$ws2.Cells.Item(3, 2).Value = 1103261.75

# Row 4: BRVM - AUTRES SECTEURS
$ws2.Cells.Item(4, 2).Value = 227998.5

# Row 5: BRVM - DISTRIBUTION
$ws2.Cells.Item(5, 2).Value = 131909.64

# Row 6: BRVM - TRANSPORT
$ws2.Cells.Item(6, 2).Value = 44122.16

# Row 7: BRVM - AGRICULTURE
$ws2.Cells.Item(7, 2).Value = 36681.95

# Row 8: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws2.Cells.Item(8, 2).Value = 5224.46

# Row 9: BRVM - FINANCES
$ws2.Cells.Item(9, 2).Value = 3629.31

# Row 10: BRVM - SERVICES FINANCIERS
$ws2.Cells.Item(10, 2).Value = 3478.92

# Row 11: BRVM-PRESTIGE
$ws2.Cells.Item(11, 2).Value = 3323.7

